$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.290.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.566.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.44%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.84%  '

$ws.Range("E6").Value = '  -0.49%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.477'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.01%  '

$ws.Range("E8").Value = '  -2.50%  '

$ws.Range("E9").Value = '  -2.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.57%  '

$ws.Range("E11").Value = '  -0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.783.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.569.39'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.42%  '

$ws.Range("E15").Value = '  -3.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.288.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '59.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.90%  '

$ws.Range("E18").Value = '  -3.29%  '

$ws.Range("E19").Value = '  -0.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '185.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.29'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.99%  '

$ws.Range("E24").Value = '  -1.94%  '

$ws.Range("E25").Value = '  -0.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.74%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.97%  '

$ws.Range("E28").Value = '  -1.91%  '

$ws.Range("E29").Value = '  -3.78%  '

$ws.Range("E30").Value = '  -6.04%  '

$ws.Range("E31").Value = '  -3.71%  '

$ws.Range("E32").Value = '  -2.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.54%  '

$ws.Range("E34").Value = '  -1.80%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.30'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.087.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.41%  '

$ws.Range("E37").Value = '  -0.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.31'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.81%  '

$ws.Range("E39").Value = '  -2.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.496'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.774'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.763'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '93.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.31%  '

$ws.Range("E44").Value = '  -2.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.696.72'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.55%  '

$ws.Range("E46").Value = '  -2.56%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '52.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.00%  '

$ws.Range("E48").Value = '  -3.55%  '

$ws.Range("E49").Value = '  -3.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.406'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.77%  '

$ws.Range("E51").Value = '  -0.61%  '
